# Generate Report for Handback
# - Flips the Status column from "Ready for handoff" to "Handed back: in sync with en-US"
# - Populates the "Latest Target File" / "Latest Handback File" columns (F/G) with
#   hyperlinked copies of the handoff md / xlf files for each locale sheet
# - Stamps the "Latest Handback DateTime" column (H) with the handback timestamp,
#   per-locale-sheet

$wb = $excel.ActiveWorkbook

$handedBackStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (touches the Overview summary sheet as well as both locale sheets)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $handedBackStatus
$wsOverview.Range("C2").Value = $handedBackStatus
$wsOverview.Range("B3").Value = $handedBackStatus
$wsOverview.Range("C3").Value = $handedBackStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $handedBackStatus
$wsZhCn.Range("C3").Value = $handedBackStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $handedBackStatus
$wsDeDe.Range("C3").Value = $handedBackStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (F) / Latest Handback File (G)
#    for both data rows, each hyperlinked the same way the original handoff
#    md / xlf links were, and stamp the handback datetime (H)
# ---------------------------------------------------------------------------
$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2f0b2b6838b566f12eac6c458978875f3fbae8fe/e2e/beae0205-6e1d-4be2-9975-2e4fbaa5253f.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b8bc734aab060b0ad35510c73671712dcc63be45/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/beae0205-6e1d-4be2-9975-2e4fbaa5253f.ccd4440b0fce08797f252197d6afa955b4655a16.zh-cn.xlf"
$zhMdName = "beae0205-6e1d-4be2-9975-2e4fbaa5253f.md"
$zhXlfName = "beae0205-6e1d-4be2-9975-2e4fbaa5253f.ccd4440b0fce08797f252197d6afa955b4655a16.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhMdUrl, "", "", $zhMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhMdUrl, "", "", $zhMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

$wsZhCn.Range("H2").Value = "2016-03-12 08:45:07"
$wsZhCn.Range("H3").Value = "2016-03-12 08:45:07"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same treatment, using the de-de handoff links, and its own
#    handback timestamp
# ---------------------------------------------------------------------------
$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2f0b2b6838b566f12eac6c458978875f3fbae8fe/e2e/beae0205-6e1d-4be2-9975-2e4fbaa5253f.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf4596711695c09006148cc20f3c006020fc3cc9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/beae0205-6e1d-4be2-9975-2e4fbaa5253f.ccd4440b0fce08797f252197d6afa955b4655a16.de-de.xlf"
$deMdName = "beae0205-6e1d-4be2-9975-2e4fbaa5253f.md"
$deXlfName = "beae0205-6e1d-4be2-9975-2e4fbaa5253f.ccd4440b0fce08797f252197d6afa955b4655a16.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deMdUrl, "", "", $deMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deXlfUrl, "", "", $deXlfName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deMdUrl, "", "", $deMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deXlfUrl, "", "", $deXlfName)

$wsDeDe.Range("H2").Value = "2016-03-12 08:45:13"
$wsDeDe.Range("H3").Value = "2016-03-12 08:45:13"
